$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''255.95'
$ws.Range("E2").Value = '''-0.69%'
$ws.Range("D3").Value = '''26.96'
$ws.Range("E3").Value = '''-0.23%'
$ws.Range("D4").Value = '''4.327'
$ws.Range("E4").Value = '''-7.77%'
$ws.Range("D5").Value = '''0.05881'
$ws.Range("E5").Value = '''-1.50%'
$ws.Range("D6").Value = '''6.619'
$ws.Range("E6").Value = '''-0.84%'
$ws.Range("D7").Value = '''0.8506'
$ws.Range("E7").Value = '''-2.37%'
$ws.Range("D8").Value = '''0.9316'
$ws.Range("E8").Value = '''-2.58%'
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = '''0.0006062'
$ws.Range("E9").Value = '''-1.00%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1379'
$ws.Range("E10").Value = '''-2.20%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.04615'
$ws.Range("E11").Value = '''27.60%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.07085'
$ws.Range("E12").Value = '''-1.43%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03064'
$ws.Range("E13").Value = '''-2.63%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09104'
$ws.Range("E14").Value = '''-1.51%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001523'
$ws.Range("E15").Value = '''-1.16%'
$ws.Range("D16").Value = '''0.006188'
$ws.Range("E16").Value = '''3.93%'
$ws.Range("D17").Value = '''3.483'
$ws.Range("E17").Value = '''-0.10%'
$ws.Range("D18").Value = '''3.169'
$ws.Range("E18").Value = '''-1.21%'
$ws.Range("E19").Value = '''-0.68%'
$ws.Range("E20").Value = '''-2.01%'
$ws.Range("D21").Value = '''0.1269'
$ws.Range("E21").Value = '''-2.87%'
$ws.Range("D22").Value = '''3.917'
$ws.Range("D23").Value = '''0.04270'
$ws.Range("E23").Value = '''0.70%'
$ws.Range("D24").Value = '''0.001218'
$ws.Range("E24").Value = '''-0.61%'
$ws.Range("D25").Value = '''0.004287'
$ws.Range("E25").Value = '''-4.83%'
$ws.Range("E26").Value = '''-0.03%'
$ws.Range("D27").Value = '''0.0001524'
$ws.Range("E27").Value = '''2.02%'
$ws.Range("D40").Value = '''0.03806'
$ws.Range("E40").Value = '''-0.54%'
$ws.Range("D41").Value = '''0.006272'
$ws.Range("E41").Value = '''56.59%'
$ws.Range("D42").Value = '''0.1100'
$ws.Range("E42").Value = '''-0.31%'
$ws.Range("D43").Value = '''0.002201'
$ws.Range("E43").Value = '''-4.37%'
$ws.Range("D44").Value = '''0.01350'
$ws.Range("E44").Value = '''22.83%'
$ws.Range("D45").Value = '''0.00005370'
$ws.Range("E45").Value = '''-2.30%'
$ws.Range("E46").Value = '''-0.02%'
$ws.Range("D47").Value = '''0.05501'
$ws.Range("E47").Value = '''-49.60%'
$ws.Range("E48").Value = '''10,965.17%'
$ws.Range("D49").Value = '''0.00002101'
$ws.Range("E49").Value = '''-0.02%'
$ws.Range("D50").Value = '''0.0002001'
$ws.Range("E50").Value = '''-0.02%'
